# Rename speaker tag "R1" to "T" in column D (Speaker column) of the active sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$usedRange = $ws.UsedRange
$lastRow = $usedRange.Rows.Count + $usedRange.Row - 1

for ($r = 1; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 4)  # Column D
    if ($cell.Value2 -eq "R1") {
        $cell.Value2 = "T"
    }
}
